# Update cryptocurrency price/volume figures scraped from coinranking.com.
# NumberFormat is forced to "@" (Text) before each Price/Volume write so
# Excel does not auto-convert numeric-looking strings (losing formatting
# like trailing zeros), matching the original inline-string cell type.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "43.553.40"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -1.49%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.289.32"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +1.24%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "95.22"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "267.63"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -2.68%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.96%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -4.12%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "44.53"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -8.20%  "
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -1.11%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.79"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -5.18%  "
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +0.27%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.630.88"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +1.23%  "
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -2.74%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.849"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +1.25%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.290.17"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +1.72%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "43.581.12"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -1.41%  "
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.02%  "
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -0.54%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "72.51"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +2.13%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.46"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +4.44%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "234.98"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -0.15%  "
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -16.68%  "
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -0.08%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.20"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -2.84%  "
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +2.12%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "40.56"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +0.88%  "
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -0.19%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "175.28"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +1.00%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "21.93"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +3.17%  "
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -4.19%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.35"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -7.12%  "
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +0.57%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0357"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -0.02%  "
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -5.86%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.38"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -1.12%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.30"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -8.33%  "
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +6.29%  "
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -7.08%  "
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +14.99%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "64.70"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +2.98%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "12.01"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -5.35%  "
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +2.51%  "
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -4.97%  "
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -1.82%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "98.11"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -2.40%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.19"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -0.80%  "

# Row 50/51: RocketPoolETH and Stacks swapped positions in the ranking,
# each carrying its own updated price/volume figures.
$ws.Range("B50").Value = "Stacks"
$ws.Range("C50").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.49"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +4.36%  "

$ws.Range("B51").Value = "RocketPoolETH"
$ws.Range("C51").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.509.99"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +1.28%  "
